$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose textual content changes in this update. Every one of these
# cells holds plain text (coin names, URLs, price strings, percentage
# strings) in the original workbook, so we explicitly force a Text number
# format before writing the new value. This stops Excel from "helpfully"
# reinterpreting strings such as "7.40" or "0.999" as numbers (which would
# silently drop significant trailing/formatting digits). We then restore
# the default "Normal" cell style so we do not leave any stray formatting
# behind that was not present in the target workbook.

$changedCells = @(
    "D2", "E2", "D3", "E3", "E4", "D5", "E5", "D6", "E6", "E7",
    "D8", "E8", "E9", "D10", "E10", "E11", "E12", "D13", "E14", "E15",
    "E16", "D17", "E17", "D18", "E18", "E19", "D20", "E20", "E21", "E22",
    "D23", "E23", "E24", "D25", "E25", "E26", "D27", "E27", "D29", "E29",
    "E30", "E31", "E32", "D33", "E33", "E34", "D35", "E35", "B36", "C36",
    "D36", "E36", "B37", "C37", "D37", "E37", "E38", "D39", "E39", "D40",
    "E40", "D41", "E41", "E42", "B43", "C43", "D43", "E43", "B44", "C44",
    "D44", "E44", "E45", "E46", "D47", "E47", "D48", "E48", "E49", "D50",
    "E50", "D51", "E51"
)

foreach ($cellRef in $changedCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "62.440.51"
$ws.Range("E2").Value = "  +4.04%  "
$ws.Range("D3").Value = "3.335.51"
$ws.Range("E3").Value = "  +4.10%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "561.45"
$ws.Range("E5").Value = "  +4.52%  "
$ws.Range("D6").Value = "151.46"
$ws.Range("E6").Value = "  +4.11%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.336.83"
$ws.Range("E8").Value = "  +3.90%  "
$ws.Range("E9").Value = "  +0.44%  "
$ws.Range("D10").Value = "7.40"
$ws.Range("E10").Value = "  +0.73%  "
$ws.Range("E11").Value = "  +3.76%  "
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").Value = "3.910.51"
$ws.Range("E14").Value = "  +0.36%  "
$ws.Range("E15").Value = "  +2.97%  "
$ws.Range("E16").Value = "  +2.65%  "
$ws.Range("D17").Value = "62.447.49"
$ws.Range("E17").Value = "  +3.94%  "
$ws.Range("D18").Value = "3.332.18"
$ws.Range("E18").Value = "  +4.38%  "
$ws.Range("E19").Value = "  +1.30%  "
$ws.Range("D20").Value = "13.77"
$ws.Range("E20").Value = "  +4.44%  "
$ws.Range("E21").Value = "  +0.67%  "
$ws.Range("E22").Value = "  +1.85%  "
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("E24").Value = "  +1.47%  "
$ws.Range("D25").Value = "69.86"
$ws.Range("E25").Value = "  -0.32%  "
$ws.Range("E26").Value = "  +5.16%  "
$ws.Range("D27").Value = "8.98"
$ws.Range("E27").Value = "  +2.53%  "
$ws.Range("D29").Value = "0.0₃0944"
$ws.Range("E29").Value = "  +4.98%  "
$ws.Range("E30").Value = "  +6.81%  "
$ws.Range("E31").Value = "  +3.85%  "
$ws.Range("E32").Value = "  +3.63%  "
$ws.Range("D33").Value = "22.85"
$ws.Range("E33").Value = "  +2.13%  "
$ws.Range("E34").Value = "  +8.21%  "
$ws.Range("D35").Value = "6.70"
$ws.Range("E35").Value = "  +0.58%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "1.47"
$ws.Range("E36").Value = "  +8.79%  "
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").Value = "159.12"
$ws.Range("E37").Value = "  +1.54%  "
$ws.Range("E38").Value = "  +12.07%  "
$ws.Range("D39").Value = "26.78"
$ws.Range("E39").Value = "  +4.52%  "
$ws.Range("D40").Value = "0.0737"
$ws.Range("E40").Value = "  +4.65%  "
$ws.Range("D41").Value = "2.790.10"
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("E42").Value = "  +7.75%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "4.25"
$ws.Range("E43").Value = "  +0.18%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").Value = "40.39"
$ws.Range("E44").Value = "  +1.13%  "
$ws.Range("E45").Value = "  +3.05%  "
$ws.Range("E46").Value = "  +4.67%  "
$ws.Range("D47").Value = "3.376.67"
$ws.Range("E47").Value = "  +4.06%  "
$ws.Range("D48").Value = "21.92"
$ws.Range("E48").Value = "  +6.17%  "
$ws.Range("E49").Value = "  -2.79%  "
$ws.Range("D50").Value = "6.29"
$ws.Range("E50").Value = "  +1.95%  "
$ws.Range("D51").Value = "286.17"
$ws.Range("E51").Value = "  +6.71%  "

foreach ($cellRef in $changedCells) {
    $ws.Range($cellRef).Style = "Normal"
}
